# edit.ps1 - apply the two changes described by the diff:
#   1. Fix the typo "tokem" -> "token" in the heading
#      "How to configure access tokem for awscli", ending up with the
#      text split across three runs: "How to configure access toke" / "n" / " for awscli"
#   2. Flip the "Normal" style's overflowPunct paragraph property from false to true.

$d = $word.ActiveDocument

# --- 1. Fix "tokem" -> "token" -------------------------------------------------

# Locate the misspelled word so we don't depend on hard-coded character offsets.
$findRng = $d.Content
$null = $findRng.Find.Execute("tokem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$wordStart = $findRng.Start   # first character of "tokem"
$mStart = $wordStart + 4      # the offending "m" is the 5th character of "tokem"
$mEnd = $wordStart + 5

# Whole containing paragraph (text only, i.e. without the trailing paragraph mark).
$para = $findRng.Paragraphs(1).Range
$paraStart = $para.Start
$paraTextEnd = $para.End - 1   # exclude the paragraph-mark character

# Replace the "m" with "n" in place (turns "tokem" into "token").
$mRng = $d.Range($mStart, $mEnd)
$mRng.Text = "n"

# The engine merges adjacent same-formatted runs back together, so nudge each of
# the three pieces' character formatting (on, then back off) to force them to stay
# as separate runs, matching the three-run split
# ("How to configure access toke" / "n" / " for awscli")
# produced by the original interactive edit.
$beforeRng = $d.Range($paraStart, $mStart)   # "How to configure access toke"
$beforeRng.Font.Bold = $true
$beforeRng.Font.Bold = $false

$nRng = $d.Range($mStart, $mStart + 1)       # "n"
$nRng.Font.Bold = $true
$nRng.Font.Bold = $false

$afterRng = $d.Range($mStart + 1, $paraTextEnd)  # " for awscli"
$afterRng.Font.Bold = $true
$afterRng.Font.Bold = $false

# --- 2. Normal style: overflowPunct false -> true ------------------------------

$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $true
